$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing row 3 (CageID 13A)
$ws.Range("B3").Value = 25
$ws.Range("C3").Value = 40
$ws.Range("D3").Value = 60
$ws.Range("E3").Value = "Plastic"

# Update existing row 5 (CageID 50A)
$ws.Range("B5").Value = 23
$ws.Range("C5").Value = 23
$ws.Range("D5").Value = 23
$ws.Range("E5").Value = "Metal"

# Update existing row 11 (CageID 15R)
$ws.Range("B11").Value = 224
$ws.Range("C11").Value = 124
$ws.Range("D11").Value = 157
$ws.Range("E11").Value = "Metal"

# Add new row 33 with cage/bird info
$ws.Range("A33").Value = 616
$ws.Range("B33").Value = 100
$ws.Range("C33").Value = 100
$ws.Range("D33").Value = 100
$ws.Range("E33").Value = "Wood"
